$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.117.96"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.834.15"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'243.47"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.07494"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'23.23"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.832.27"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'0.6673"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'0.000009367"
$ws.Range("E16").Value = "  -7.80%  "
$ws.Range("D17").Value = "'5.983"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "29.115.76"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "2.081.67"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'223.13"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'7.099"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'160.16"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "'8.489"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'17.88"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "'0.05649"
$ws.Range("E30").Value = "  +7.70%  "
$ws.Range("D31").Value = "'4.153"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "'4.082"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "'0.7418"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'1.836"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'1.141"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'2.673"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'2.760"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "1.219.28"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").Value = "'0.01778"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "'0.8922"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'101.95"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.982.30"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'65.67"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "'0.00000000123"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "'0.5096"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "'0.4071"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "'0.07416"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").Value = "'9.002"
$ws.Range("E51").Value = "  +1.63%  "
